# "Generate Report for Handback" — records a handback-report result for the
# fd3e74e1-... row on both the zh-cn and de-de sheets:
#   - "Latest Target File" (col I) becomes a hyperlink to the handback .md file
#   - "Latest Handback File" (col J) is stamped with the same value as
#     "Latest Target File" (col H)
#   - "Latest Handback DateTime" (col K) gets a fresh timestamp
#   - "Error Detail" (col P) records that the handback file version is stale
# Columns I, J and P are also widened to match the other wide columns (40).

$wb = $excel.ActiveWorkbook

$mdDisplay = "fd3e74e1-b6c4-4fa1-9a20-3da042eef3f2.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ebfa86a689a9e1b4f630af7a09a3c79e8d5453/e2e/fd3e74e1-b6c4-4fa1-9a20-3da042eef3f2.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f20fbe22e3821466bad32cefe11fd9332e4db404/e2e/fd3e74e1-b6c4-4fa1-9a20-3da042eef3f2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ebfa86a689a9e1b4f630af7a09a3c79e8d5453/e2e/fd3e74e1-b6c4-4fa1-9a20-3da042eef3f2.md."

# Hyperlink font (matches the workbook's existing "HyperLink" cell style:
# underlined Calibri 11 in FF6495ED).
$hyperlinkUnderline = 2        # xlUnderlineStyleSingle
$hyperlinkColor     = 15570276 # RGB(0x64, 0x95, 0xED) -> matches FF6495ED

function Set-HandbackReport {
    param(
        [string]$SheetName,
        [string]$KDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen columns I (9), J (10) and P (16) to 40 characters, same as the
    # other "wide" columns on this sheet.
    $ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
    $ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(1).ColumnWidth
    $ws.Columns.Item(16).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

    # I6: "Latest Target File" -> hyperlink to the handback markdown file.
    $i6 = $ws.Cells.Item(6, 9)
    $ws.Hyperlinks.Add($i6, $mdUrl, "", "", $mdDisplay) | Out-Null
    $i6.Font.Underline = $hyperlinkUnderline
    $i6.Font.Color = $hyperlinkColor

    # J6: "Latest Handback File" -> same value as "Latest Target File" (H6).
    $ws.Cells.Item(6, 10).Value = $ws.Cells.Item(6, 8).Value

    # K6: "Latest Handback DateTime" -> new timestamp.
    $ws.Cells.Item(6, 11).Value = $KDateTime

    # P6: "Error Detail" -> stale-handback-version message.
    $ws.Cells.Item(6, 16).Value = $errorDetail
}

Set-HandbackReport "zh-cn" "2016-10-27 07:53:35"
Set-HandbackReport "de-de" "2016-10-27 07:53:52"
